# Build site at 2022-09-26 16:07:08 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (column B / C data cells) ---
$ws.Range("B10").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C10").Value = "5840560 - Marco Antonio Carvalho Pereira"

$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = ""

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"

$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = ""

$ws.Range("A17").Value = "Avaliação:"

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C18").Value = "5840560 - Marco Antonio Carvalho Pereira"

$ws.Range("A19").Value = "Critério:"

$ws.Range("A20").Value = "Norma de recuperação:"

$ws.Range("A21").Value = "Bibliografia:"

$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = ""

$ws.Range("A23").Value = ""
$ws.Range("B23").Value = "LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)`n"

# --- Delete the now-superfluous last row (old row 24) ---
$ws.Rows.Item(24).Delete()

# --- Fix up row heights to match the new layout ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 15
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).RowHeight = 15
$ws.Rows.Item(23).RowHeight = 30
